# Apply scheduled runner price/profit updates to Cerberus_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 250000000
$ws.Range("I62").Value = 250000000
$ws.Range("K62").Value = 250000000
$ws.Range("M62").Value = -249999376

# Row 65
$ws.Range("H65").Value = 250000000
$ws.Range("I65").Value = 250000000
$ws.Range("K65").Value = 1250000000
$ws.Range("M65").Value = -1249996880

# Row 100
$ws.Range("H100").Value = 16910.6
$ws.Range("J100").Value = 16910.6
$ws.Range("L100").Value = 16910.6
$ws.Range("N100").Value = -17992.6

# Row 103
$ws.Range("H103").Value = 980.9091
$ws.Range("J103").Value = 998.8889
$ws.Range("L103").Value = 2996.6667
$ws.Range("N103").Value = -4168.6667

# Row 111
$ws.Range("H111").Value = 2422.1538
$ws.Range("I111").Value = 2405.7
$ws.Range("J111").Value = 2477
$ws.Range("K111").Value = 7217.099999999999
$ws.Range("L111").Value = 7431
$ws.Range("M111").Value = -4150.099999999999
$ws.Range("N111").Value = -13565

# Row 125
$ws.Range("H125").Value = 2453.4375
$ws.Range("I125").Value = 2755.2856
$ws.Range("J125").Value = 2218.6667
$ws.Range("K125").Value = 24797.5704
$ws.Range("L125").Value = 19968.0003
$ws.Range("M125").Value = -22337.5704
$ws.Range("N125").Value = -24888.0003

# Row 132
$ws.Range("H132").Value = 3187.8333
$ws.Range("I132").Value = 3036.6274
$ws.Range("K132").Value = 9109.8822
$ws.Range("M132").Value = -6579.8822

# Row 134
$ws.Range("H134").Value = 33998.785
$ws.Range("J134").Value = 33998.785
$ws.Range("L134").Value = 33998.785
$ws.Range("N134").Value = -44138.785

# Row 137
$ws.Range("H137").Value = 2549.682
$ws.Range("J137").Value = 3012
$ws.Range("L137").Value = 9036
$ws.Range("N137").Value = -14136

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2691.4
$ws.Range("I45").Value = 2198.2222
$ws.Range("K45").Value = 2198.2222
$ws.Range("M45").Value = -1821.2222

# Row 61
$ws.Range("H61").Value = 4235.154
$ws.Range("I61").Value = 3595.2222
$ws.Range("K61").Value = 3595.2222
$ws.Range("M61").Value = -3383.2222

# Row 63
$ws.Range("H63").Value = 2876.4
$ws.Range("I63").Value = 2640.4443
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2640.4443
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1954.4443
$ws.Range("N63").Value = -6372

# Row 66
$ws.Range("H66").Value = 2876.4
$ws.Range("I66").Value = 2640.4443
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 13202.2215
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -9770.2215
$ws.Range("N66").Value = -31864

# Row 88
$ws.Range("H88").Value = 9125.615
$ws.Range("I88").Value = 1376.7142
$ws.Range("K88").Value = 1376.7142
$ws.Range("M88").Value = -970.7141999999999

# Row 91
$ws.Range("H91").Value = 9125.615
$ws.Range("I91").Value = 1376.7142
$ws.Range("K91").Value = 1376.7142
$ws.Range("M91").Value = 27.28580000000011

# Row 136
$ws.Range("H136").Value = 4235.154
$ws.Range("I136").Value = 3595.2222
$ws.Range("K136").Value = 10785.6666
$ws.Range("M136").Value = -8235.6666

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 13191.272
$ws.Range("I86").Value = 3888.125
$ws.Range("J86").Value = 37999.668
$ws.Range("K86").Value = 3888.125
$ws.Range("L86").Value = 37999.668
$ws.Range("M86").Value = -2765.125
$ws.Range("N86").Value = -40245.668

# Row 89
$ws.Range("H89").Value = 13191.272
$ws.Range("I89").Value = 3888.125
$ws.Range("J89").Value = 37999.668
$ws.Range("K89").Value = 19440.625
$ws.Range("L89").Value = 189998.34
$ws.Range("M89").Value = -13824.625
$ws.Range("N89").Value = -201230.34

# Row 134
$ws.Range("H134").Value = 10222.192
$ws.Range("I134").Value = 8964
$ws.Range("J134").Value = 17142.25
$ws.Range("K134").Value = 26892
$ws.Range("L134").Value = 51426.75
$ws.Range("M134").Value = -24357
$ws.Range("N134").Value = -56496.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3154.375
$ws.Range("I31").Value = 2932.375
$ws.Range("J31").Value = 3820.375
$ws.Range("K31").Value = 2932.375
$ws.Range("L31").Value = 3820.375
$ws.Range("M31").Value = -2637.375
$ws.Range("N31").Value = -4410.375

# Row 34
$ws.Range("H34").Value = 3154.375
$ws.Range("I34").Value = 2932.375
$ws.Range("J34").Value = 3820.375
$ws.Range("K34").Value = 2932.375
$ws.Range("L34").Value = 3820.375
$ws.Range("M34").Value = -2730.375
$ws.Range("N34").Value = -4224.375

# Row 62
$ws.Range("H62").Value = 29510.25
$ws.Range("J62").Value = 29510.25
$ws.Range("L62").Value = 29510.25
$ws.Range("N62").Value = -30758.25

# Row 65
$ws.Range("H65").Value = 29510.25
$ws.Range("J65").Value = 29510.25
$ws.Range("L65").Value = 147551.25
$ws.Range("N65").Value = -153791.25

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 2825042.2
$ws.Range("I11").Value = 3390028.5
$ws.Range("J11").Value = 111
$ws.Range("K11").Value = 10170085.5
$ws.Range("L11").Value = 333
$ws.Range("M11").Value = -10169945.5
$ws.Range("N11").Value = -613

# Row 22
$ws.Range("H22").Value = 1442.25
$ws.Range("J22").Value = 1916.3334
$ws.Range("L22").Value = 5749.0002
$ws.Range("N22").Value = -6087.0002

# Row 27
$ws.Range("H27").Value = 1442.25
$ws.Range("J27").Value = 1916.3334
$ws.Range("L27").Value = 5749.0002
$ws.Range("N27").Value = -5953.0002

# Row 137
$ws.Range("H137").Value = 17149.25
$ws.Range("I137").Value = 18148
$ws.Range("K137").Value = 54444
$ws.Range("M137").Value = -49344

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 1766
$ws.Range("I80").Value = 2054
$ws.Range("J80").Value = 1382
$ws.Range("K80").Value = 2054
$ws.Range("L80").Value = 1382
$ws.Range("M80").Value = -1056
$ws.Range("N80").Value = -3378

# Row 83
$ws.Range("H83").Value = 1766
$ws.Range("I83").Value = 2054
$ws.Range("J83").Value = 1382
$ws.Range("K83").Value = 10270
$ws.Range("L83").Value = 6910
$ws.Range("M83").Value = -5278
$ws.Range("N83").Value = -16894

# Row 135
$ws.Range("H135").Value = 87809.60000000001
$ws.Range("I135").Value = 50709
$ws.Range("J135").Value = 97084.75
$ws.Range("K135").Value = 50709
$ws.Range("L135").Value = 97084.75
$ws.Range("M135").Value = -45639
$ws.Range("N135").Value = -107224.75

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2868.9
$ws.Range("I61").Value = 2742.625
$ws.Range("K61").Value = 2742.625
$ws.Range("M61").Value = -2540.625

# Row 93
$ws.Range("H93").Value = 1611.8572
$ws.Range("I93").Value = 1096.8
$ws.Range("K93").Value = 1096.8
$ws.Range("M93").Value = 151.2

# Row 100
$ws.Range("H100").Value = 2492
$ws.Range("I100").Value = 1427.1428
$ws.Range("K100").Value = 1427.1428
$ws.Range("M100").Value = -886.1428000000001

# Row 113
$ws.Range("H113").Value = 2868.9
$ws.Range("I113").Value = 2742.625
$ws.Range("K113").Value = 2742.625
$ws.Range("M113").Value = -572.625

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8501
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 8501
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 8501
$ws.Range("N62").Value = -9749
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 8501
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 8501
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 42505
$ws.Range("N65").Value = -48745
$ws.Range("M65").ClearContents()

# Row 81
$ws.Range("H81").Value = 7852.875
$ws.Range("I81").Value = 7852.875
$ws.Range("K81").Value = 15705.75
$ws.Range("M81").Value = -14644.75

# Row 84
$ws.Range("H84").Value = 7852.875
$ws.Range("I84").Value = 7852.875
$ws.Range("K84").Value = 78528.75
$ws.Range("M84").Value = -73224.75

# Row 107
$ws.Range("H107").Value = 1145.6666
$ws.Range("I107").Value = 462.4
$ws.Range("K107").Value = 1387.2
$ws.Range("M107").Value = 532.8000000000002
